# ---------------------------------------------------------------------------
# Adds a new worksheet "0.9" (an ER-diagram attribute/table listing) in front
# of the existing "Sheet1", which is renamed to "v0.5".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet to "v0.5" and capture its old selection ----
$v05 = $wb.Worksheets.Item(1)
$v05.Activate()
$v05.Range("B26:H27").Select()
$v05.Name = "v0.5"

# --- Create the new front sheet "0.9" --------------------------------------
$new = $wb.Worksheets.Add()
$new.Name = "0.9"

# Page setup to match the rest of the workbook
$new.PageSetup.LeftMargin = 54
$new.PageSetup.RightMargin = 54
$new.PageSetup.TopMargin = 72
$new.PageSetup.BottomMargin = 72
$new.PageSetup.HeaderMargin = 36
$new.PageSetup.FooterMargin = 36
$new.PageSetup.Orientation = 1

# Column widths (character units; engine stores width+5/6)
$new.Columns.Item(2).ColumnWidth = 12.666666666666666   # B -> 13.5
$new.Columns.Item(3).ColumnWidth = 9.666666666666666    # C -> 10.5
$new.Columns.Item(4).ColumnWidth = 10.330729166666666   # D -> ~11.164
$new.Columns.Item(5).ColumnWidth = 12.330729166666666   # E -> ~13.164
$new.Columns.Item(6).ColumnWidth = 7.666666666666667    # F -> 8.5
$new.Columns.Item(8).ColumnWidth = 7.330729166666667    # H -> ~8.164

# ===========================================================================
# Table: Employee
# ===========================================================================
$new.Range("B2:H2").Merge()
$new.Range("B2").Value = "Employee"
$new.Range("B2:H2").Font.Bold = $true
$new.Range("B2:H2").HorizontalAlignment = -4108

$new.Range("B3").Value = "Fname"
$new.Range("C3").Value = "Minit"
$new.Range("D3").Value = "Lname"
$new.Range("E3").Value = "Terminated"
$new.Range("F3").Value = "E-ID"
$new.Range("G3").Value = "Manager ID"
$new.Range("H3").Value = "hash pin"
$new.Range("I3").Value = "Derpartment"

# Row 4 stays blank (kept for spacing, same row-level formatting as row 3)
$new.Range("B4:H4").Value = ""

# ===========================================================================
# Table: Mesagge
# ===========================================================================
$new.Range("B7:H7").Merge()
$new.Range("B7").Value = "Mesagge"
$new.Range("B7:H7").Font.Bold = $true
$new.Range("B7:H7").HorizontalAlignment = -4108

$new.Range("B8").Value = "M-ID"
$new.Range("C8").Value = "message"
$new.Range("D8:E8").Merge()
$new.Range("D8").Value = "Msg_from_E-ID"
$new.Range("D8:E8").HorizontalAlignment = -4108

# ===========================================================================
# Table: Mesagges_For
# ===========================================================================
$new.Range("B11:H11").Merge()
$new.Range("B11").Value = "Mesagges_For"
$new.Range("B11:H11").Font.Bold = $true
$new.Range("B11:H11").HorizontalAlignment = -4108

$new.Range("B12").Value = "M-ID"
$new.Range("C12").Value = "E-ID"
$new.Range("D12").Value = "Pending"

# ===========================================================================
# Table: Department
# ===========================================================================
$new.Range("B16:H16").Merge()
$new.Range("B16").Value = "Department"
$new.Range("B16:H16").Font.Bold = $true
$new.Range("B16:H16").HorizontalAlignment = -4108

$new.Range("B17").Value = "D-ID"
$new.Range("C17").Value = "D-Name"
$new.Range("D17").Value = "Location"
$new.Range("E17").Value = "More?"
$new.Range("E17").Font.Color = 255

# ===========================================================================
# Table: Holiday
# ===========================================================================
$new.Range("B21:H21").Merge()
$new.Range("B21").Value = "Holiday"
$new.Range("B21:H21").Font.Bold = $true
$new.Range("B21:H21").HorizontalAlignment = -4108

$new.Range("B22").Value = "H-ID"
$new.Range("C22").Value = "Date"
$new.Range("D22").Value = "Repeats"
$new.Range("E22").Value = "Not totally sure how to make this in a good way"
$new.Range("E22").Font.Color = 255

# ===========================================================================
# Table: Hoplidays_held
# ===========================================================================
$new.Range("B25:H25").Merge()
$new.Range("B25").Value = "Hoplidays_held"
$new.Range("B25:H25").Font.Bold = $true
$new.Range("B25:H25").HorizontalAlignment = -4108

$new.Range("B26").Value = "D-ID"
$new.Range("C26").Value = "H-ID"

# ===========================================================================
# Table: Pay_Type
# ===========================================================================
$new.Range("B30:H30").Merge()
$new.Range("B30").Value = "Pay_Type"
$new.Range("B30:H30").Font.Bold = $true
$new.Range("B30:H30").HorizontalAlignment = -4108

$new.Range("B31").Value = "Pay_Type-ID"
$new.Range("C31").Value = "Daily max"
$new.Range("D31").Value = "Weekly max"
$new.Range("E31:F31").Merge()
$new.Range("E31").Value = "Pay_type_when over time"
$new.Range("E31:F31").HorizontalAlignment = -4108

# Lone formatted (but empty) cell between tables
$new.Range("B34").Value = ""
$new.Range("B34").Font.Bold = $true

# ===========================================================================
# Table: Pay_Type in Departments
# ===========================================================================
$new.Range("B35:H35").Merge()
$new.Range("B35").Value = "Pay_Type in Departments"
$new.Range("B35:H35").Font.Bold = $true
$new.Range("B35:H35").HorizontalAlignment = -4108

$new.Range("B36").Value = "Pay_Type-ID"
$new.Range("C36").Value = "D-ID"

# ===========================================================================
# Table: Punch
# ===========================================================================
$new.Range("B40:H40").Merge()
$new.Range("B40").Value = "Punch"
$new.Range("B40:H40").Font.Bold = $true
$new.Range("B40:H40").HorizontalAlignment = -4108

$new.Range("B41").Value = "P-ID"
$new.Range("C41").Value = "In-time"
$new.Range("D41").Value = "out-time"
$new.Range("E41").Value = "Punch_type-ID"
$new.Range("F41").Value = "D-ID"
$new.Range("G41").Value = "E-ID"

# ===========================================================================
# Table: Punch_Type
# ===========================================================================
$new.Range("B45:G45").Merge()
$new.Range("B45").Value = "Punch_Type"
$new.Range("B45:G45").Font.Bold = $true
$new.Range("B45:G45").HorizontalAlignment = -4108

$new.Range("B46").Value = "Punch_Type-ID"
$new.Range("C46").Value = "Description"
$new.Range("D46:E46").Merge()
$new.Range("D46").Value = "punch_in_option"
$new.Range("D46:E46").HorizontalAlignment = -4108

$new.Range("D47:E47").HorizontalAlignment = -4108
$new.Range("D48:E48").HorizontalAlignment = -4108

# ===========================================================================
# Table: Timecard
# ===========================================================================
$new.Range("B50:H50").Merge()
$new.Range("B50").Value = "Timecard"
$new.Range("B50:H50").Font.Bold = $true
$new.Range("B50:H50").HorizontalAlignment = -4108

$new.Range("B51").Value = "E-ID"
$new.Range("C51").Value = "Pay-Period"

# ===========================================================================
# Table: Timecard Lines
# ===========================================================================
$new.Range("B55:H55").Merge()
$new.Range("B55").Value = "Timecard Lines"
$new.Range("B55:H55").Font.Bold = $true
$new.Range("B55:H55").HorizontalAlignment = -4108

$new.Range("B56").Value = "Line-Number"
$new.Range("C56").Value = "P-ID"
$new.Range("D56").Value = "E-ID"
$new.Range("E56").Value = "Pay-Period"
$new.Range("F56").Value = "Pay_Type-ID"
$new.Range("G56").Value = "Split-Start"
$new.Range("H56").Value = "Split-End"

$new.Range("I4").Select()

Write-Output "setup done"
